$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add column D "success" ---
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column D data values (text "0" / "1"), no fill-down formatting carried ---
$ws.Range("D2:D15").Value = "'0"
$ws.Range("D2:D15").Style = "Normal"
$ws.Range("D18").Value = "'0"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Value = "'0"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'0"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Value = "'0"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31:D33").Value = "'0"
$ws.Range("D31:D33").Style = "Normal"
$ws.Range("D36:D43").Value = "'0"
$ws.Range("D36:D43").Style = "Normal"
$ws.Range("D50").Value = "'0"
$ws.Range("D50").Style = "Normal"
$ws.Range("D53").Value = "'0"
$ws.Range("D53").Style = "Normal"
$ws.Range("D56:D57").Value = "'0"
$ws.Range("D56:D57").Style = "Normal"
$ws.Range("D60:D65").Value = "'0"
$ws.Range("D60:D65").Style = "Normal"
$ws.Range("D16:D17").Value = "'1"
$ws.Range("D16:D17").Style = "Normal"
$ws.Range("D19:D20").Value = "'1"
$ws.Range("D19:D20").Style = "Normal"
$ws.Range("D22").Value = "'1"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24:D26").Value = "'1"
$ws.Range("D24:D26").Style = "Normal"
$ws.Range("D28:D30").Value = "'1"
$ws.Range("D28:D30").Style = "Normal"
$ws.Range("D34:D35").Value = "'1"
$ws.Range("D34:D35").Style = "Normal"
$ws.Range("D44:D49").Value = "'1"
$ws.Range("D44:D49").Style = "Normal"
$ws.Range("D51:D52").Value = "'1"
$ws.Range("D51:D52").Style = "Normal"
$ws.Range("D54:D55").Value = "'1"
$ws.Range("D54:D55").Style = "Normal"
$ws.Range("D58:D59").Value = "'1"
$ws.Range("D58:D59").Style = "Normal"
